$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.03225138200131
$ws.Range("C2").Value = 4.482561356785101
$ws.Range("D2").Value = 7.15792718459619
$ws.Range("F2").Value = 37.47940862859942
$ws.Range("G2").Value = 3.694587258480427
$ws.Range("K2").Value = 11.4772871337697

$ws.Range("B3").Value = 11.83675915855333
$ws.Range("C3").Value = 4.328504866439691
$ws.Range("D3").Value = 7.154367573783979
$ws.Range("F3").Value = 37.17174655258061
$ws.Range("G3").Value = 3.697492045173542
$ws.Range("K3").Value = 11.36047978666874

$ws.Range("B4").Value = 11.71900950269414
$ws.Range("C4").Value = 4.232835808763559
$ws.Range("D4").Value = 7.151950977015856
$ws.Range("F4").Value = 36.98836362880832
$ws.Range("G4").Value = 3.699367698216931
$ws.Range("K4").Value = 11.29166155875524

$ws.Range("B5").Value = 11.67167023316293
$ws.Range("C5").Value = 4.193654857782745
$ws.Range("D5").Value = 7.150907660271398
$ws.Range("F5").Value = 36.91507507547308
$ws.Range("G5").Value = 3.700155283970505
$ws.Range("K5").Value = 11.26438049911395

$ws.Range("B6").Value = 11.66385056281122
$ws.Range("C6").Value = 4.187139441057168
$ws.Range("D6").Value = 7.150730869704082
$ws.Range("F6").Value = 36.9029941094659
$ws.Range("G6").Value = 3.700287468149636
$ws.Range("K6").Value = 11.25989751740477

$ws.Range("B7").Value = 11.71836836601773
$ws.Range("C7").Value = 4.232308085916158
$ws.Range("D7").Value = 7.151937144095901
$ws.Range("F7").Value = 36.98736932998906
$ws.Range("G7").Value = 3.699378225660662
$ws.Range("K7").Value = 11.29129050513202

$ws.Range("B8").Value = 11.96441427168415
$ws.Range("C8").Value = 4.429710182644111
$ws.Range("D8").Value = 7.156747500928437
$ws.Range("F8").Value = 37.37221299430865
$ws.Range("G8").Value = 3.69556976259793
$ws.Range("K8").Value = 11.4364300620723

$ws.Range("B9").Value = 12.46190536712846
$ws.Range("C9").Value = 4.805221218891453
$ws.Range("D9").Value = 7.164367597783221
$ws.Range("F9").Value = 38.16826505282251
$ws.Range("G9").Value = 3.68882842029537
$ws.Range("K9").Value = 11.74254223961419

$ws.Range("B10").Value = 12.83239244497279
$ws.Range("C10").Value = 5.070393826256428
$ws.Range("D10").Value = 7.168887211036988
$ws.Range("F10").Value = 38.77490886631035
$ws.Range("G10").Value = 3.684313527557756
$ws.Range("K10").Value = 11.97841259003044

$ws.Range("B11").Value = 13.00115805169017
$ws.Range("C11").Value = 5.188051146408537
$ws.Range("D11").Value = 7.170714585253955
$ws.Range("F11").Value = 39.05484429879479
$ws.Range("G11").Value = 3.682353566557161
$ws.Range("K11").Value = 12.08762895329144

$ws.Range("B12").Value = 13.06502910285687
$ws.Range("C12").Value = 5.232131403295763
$ws.Range("D12").Value = 7.171374143143227
$ws.Range("F12").Value = 39.16135216870574
$ws.Range("G12").Value = 3.681624795552578
$ws.Range("K12").Value = 12.12922240789226

$ws.Range("B13").Value = 13.05127605007083
$ws.Range("C13").Value = 5.222659700680139
$ws.Range("D13").Value = 7.171233530760903
$ws.Range("F13").Value = 39.13839254893207
$ws.Range("G13").Value = 3.681781153627854
$ws.Range("K13").Value = 12.12025467404799

$ws.Range("B14").Value = 13.00641388990015
$ws.Range("C14").Value = 5.191687404366334
$ws.Range("D14").Value = 7.170769494858837
$ws.Range("F14").Value = 39.06359709207351
$ws.Range("G14").Value = 3.682293341538017
$ws.Range("K14").Value = 12.09104638010701

$ws.Range("B15").Value = 12.97892769940788
$ws.Range("C15").Value = 5.1726529519285
$ws.Range("D15").Value = 7.170481049446178
$ws.Range("F15").Value = 39.01784609362291
$ws.Range("G15").Value = 3.682608817302274
$ws.Range("K15").Value = 12.07318494614725

$ws.Range("B16").Value = 12.82136180480815
$ws.Range("C16").Value = 5.062640734302732
$ws.Range("D16").Value = 7.168763225615243
$ws.Range("F16").Value = 38.75668848250599
$ws.Range("G16").Value = 3.684443498255408
$ws.Range("K16").Value = 11.97131019792706

$ws.Range("B17").Value = 12.72470849453452
$ws.Range("C17").Value = 4.994355808697718
$ws.Range("D17").Value = 7.167651159526067
$ws.Range("F17").Value = 38.59744516857064
$ws.Range("G17").Value = 3.685593007176828
$ws.Range("K17").Value = 11.90927592740455

$ws.Range("B18").Value = 12.66914166580931
$ws.Range("C18").Value = 4.954802914515452
$ws.Range("D18").Value = 7.166989983169304
$ws.Range("F18").Value = 38.5062304858558
$ws.Range("G18").Value = 3.686263015914463
$ws.Range("K18").Value = 11.87377837520612

$ws.Range("B19").Value = 12.65033425616823
$ws.Range("C19").Value = 4.941364888664347
$ws.Range("D19").Value = 7.166762406135106
$ws.Range("F19").Value = 38.47541373042944
$ws.Range("G19").Value = 3.686491390269595
$ws.Range("K19").Value = 11.86179215543796

$ws.Range("B20").Value = 12.73499526861462
$ws.Range("C20").Value = 5.001653902172254
$ws.Range("D20").Value = 7.167771767495049
$ws.Range("F20").Value = 38.61435827756574
$ws.Range("G20").Value = 3.685469725478851
$ws.Range("K20").Value = 11.91586094126471

$ws.Range("B21").Value = 13.01959253633795
$ws.Range("C21").Value = 5.20079792920513
$ws.Range("D21").Value = 7.170906670022992
$ws.Range("F21").Value = 39.08555324619529
$ws.Range("G21").Value = 3.682142535829654
$ws.Range("K21").Value = 12.099619490006

$ws.Range("B22").Value = 13.20534822140285
$ws.Range("C22").Value = 5.328169183341439
$ws.Range("D22").Value = 7.172766693860531
$ws.Range("F22").Value = 39.39640243588612
$ws.Range("G22").Value = 3.680046233135797
$ws.Range("K22").Value = 12.22107099943315

$ws.Range("B23").Value = 13.10625170548556
$ws.Range("C23").Value = 5.26045726924198
$ws.Range("D23").Value = 7.171791094076739
$ws.Range("F23").Value = 39.23025422560671
$ws.Range("G23").Value = 3.681157938640873
$ws.Range("K23").Value = 12.15613932727384

$ws.Range("B24").Value = 12.7303446113156
$ws.Range("C24").Value = 4.998355349304945
$ws.Range("D24").Value = 7.167717308668356
$ws.Range("F24").Value = 38.60671080783803
$ws.Range("G24").Value = 3.68552543265576
$ws.Range("K24").Value = 11.91288333395172

$ws.Range("B25").Value = 12.32614858671574
$ws.Range("C25").Value = 4.705276992164057
$ws.Range("D25").Value = 7.16249751600673
$ws.Range("F25").Value = 37.9488288561486
$ws.Range("G25").Value = 3.690574841759788
$ws.Range("K25").Value = 11.65764653851694
